# Insert a new weekly record for "Macroferia Regional de Talca" / Chirimoya.
# A new row is inserted at row 26 (pushing the existing rows 26-68 down to
# 27-69, carrying their data with them unchanged), then the newly opened
# row 26 is populated with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(26).Insert()

$ws.Range("A26").Value = 5
$ws.Range("B26").Value = "Macroferia Regional de Talca"
$ws.Range("C26").Value = "Maule"
$ws.Range("D26").Value = 44799
$ws.Range("E26").Value = 7
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100107
$ws.Range("H26").Value = "Otros"
$ws.Range("I26").Value = 100107002
$ws.Range("J26").Value = "Chirimoya"
$ws.Range("K26").Value = "Cultivar IV Región"
$ws.Range("L26").Value = "Primera"
$ws.Range("M26").Value = 50
$ws.Range("N26").Value = 30000
$ws.Range("O26").Value = 30000
$ws.Range("P26").Value = 30000
$ws.Range("Q26").Value = "`$/bandeja 10 kilos"
$ws.Range("R26").Value = "Provincia de Limarí"
$ws.Range("S26").Value = 3000
$ws.Range("T26").Value = 10
